$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.Value = "'26.878.23"
$c.Style = 'Normal'
$c = $ws.Range('E2')
$c.Value = "'  -1.65%  "
$c.Style = 'Normal'
$c = $ws.Range('D3')
$c.Value = "'1.826.21"
$c.Style = 'Normal'
$c = $ws.Range('E3')
$c.Value = "'  -1.63%  "
$c.Style = 'Normal'
$c = $ws.Range('D4')
$c.Value = "'1.006"
$c.Style = 'Normal'
$c = $ws.Range('E4')
$c.Value = "'  +0.52%  "
$c.Style = 'Normal'
$c = $ws.Range('D5')
$c.Value = "'310.74"
$c.Style = 'Normal'
$c = $ws.Range('E5')
$c.Value = "'  -0.95%  "
$c.Style = 'Normal'
$c = $ws.Range('D6')
$c.Value = "'1.006"
$c.Style = 'Normal'
$c = $ws.Range('E6')
$c.Value = "'  +0.52%  "
$c.Style = 'Normal'
$c = $ws.Range('D7')
$c.Value = "'0.4574"
$c.Style = 'Normal'
$c = $ws.Range('E7')
$c.Value = "'  -0.77%  "
$c.Style = 'Normal'
$c = $ws.Range('D8')
$c.Value = "'0.3678"
$c.Style = 'Normal'
$c = $ws.Range('E8')
$c.Value = "'  -0.83%  "
$c.Style = 'Normal'
$c = $ws.Range('D9')
$c.Value = "'0.07156"
$c.Style = 'Normal'
$c = $ws.Range('E9')
$c.Value = "'  -2.26%  "
$c.Style = 'Normal'
$c = $ws.Range('D10')
$c.Value = "'0.8728"
$c.Style = 'Normal'
$c = $ws.Range('E10')
$c.Value = "'  -0.89%  "
$c.Style = 'Normal'
$c = $ws.Range('D11')
$c.Value = "'0.07790"
$c.Style = 'Normal'
$c = $ws.Range('E11')
$c.Value = "'  -0.20%  "
$c.Style = 'Normal'
$c = $ws.Range('D12')
$c.Value = "'19.55"
$c.Style = 'Normal'
$c = $ws.Range('E12')
$c.Value = "'  -1.59%  "
$c.Style = 'Normal'
$c = $ws.Range('D13')
$c.Value = "'1.822.71"
$c.Style = 'Normal'
$c = $ws.Range('E13')
$c.Value = "'  -5.31%  "
$c.Style = 'Normal'
$c = $ws.Range('D14')
$c.Value = "'5.315"
$c.Style = 'Normal'
$c = $ws.Range('E14')
$c.Value = "'  -1.37%  "
$c.Style = 'Normal'
$c = $ws.Range('D15')
$c.Value = "'6.379"
$c.Style = 'Normal'
$c = $ws.Range('E15')
$c.Value = "'  -2.51%  "
$c.Style = 'Normal'
$c = $ws.Range('D16')
$c.Value = "'86.80"
$c.Style = 'Normal'
$c = $ws.Range('E16')
$c.Value = "'  -5.46%  "
$c.Style = 'Normal'
$c = $ws.Range('D17')
$c.Value = "'1.006"
$c.Style = 'Normal'
$c = $ws.Range('E17')
$c.Value = "'  +0.46%  "
$c.Style = 'Normal'
$c = $ws.Range('D18')
$c.Value = "'0.000008695"
$c.Style = 'Normal'
$c = $ws.Range('E18')
$c.Value = "'  -4.03%  "
$c.Style = 'Normal'
$c = $ws.Range('E19')
$c.Value = "'  +0.45%  "
$c.Style = 'Normal'
$c = $ws.Range('D20')
$c.Value = "'26.893.11"
$c.Style = 'Normal'
$c = $ws.Range('E20')
$c.Value = "'  -1.69%  "
$c.Style = 'Normal'
$c = $ws.Range('D21')
$c.Value = "'14.44"
$c.Style = 'Normal'
$c = $ws.Range('E21')
$c.Value = "'  -2.27%  "
$c.Style = 'Normal'
$c = $ws.Range('D22')
$c.Value = "'4.990"
$c.Style = 'Normal'
$c = $ws.Range('E22')
$c.Value = "'  -2.72%  "
$c.Style = 'Normal'
$c = $ws.Range('D23')
$c.Value = "'2.044.83"
$c.Style = 'Normal'
$c = $ws.Range('E23')
$c.Value = "'  -4.81%  "
$c.Style = 'Normal'
$c = $ws.Range('D24')
$c.Value = "'10.44"
$c.Style = 'Normal'
$c = $ws.Range('E24')
$c.Value = "'  -0.68%  "
$c.Style = 'Normal'
$c = $ws.Range('D25')
$c.Value = "'2.005"
$c.Style = 'Normal'
$c = $ws.Range('E25')
$c.Value = "'  +3.41%  "
$c.Style = 'Normal'
$c = $ws.Range('D26')
$c.Value = "'151.07"
$c.Style = 'Normal'
$c = $ws.Range('E26')
$c.Value = "'  -0.63%  "
$c.Style = 'Normal'
$c = $ws.Range('D27')
$c.Value = "'18.18"
$c.Style = 'Normal'
$c = $ws.Range('D28')
$c.Value = "'1.950"
$c.Style = 'Normal'
$c = $ws.Range('E28')
$c.Value = "'  -5.77%  "
$c.Style = 'Normal'
$c = $ws.Range('D29')
$c.Value = "'113.47"
$c.Style = 'Normal'
$c = $ws.Range('D30')
$c.Value = "'4.909"
$c.Style = 'Normal'
$c = $ws.Range('E30')
$c.Value = "'  -3.81%  "
$c.Style = 'Normal'
$c = $ws.Range('D31')
$c.Value = "'0.08794"
$c.Style = 'Normal'
$c = $ws.Range('E31')
$c.Value = "'  -0.78%  "
$c.Style = 'Normal'
$c = $ws.Range('D32')
$c.Value = "'3.017"
$c.Style = 'Normal'
$c = $ws.Range('E32')
$c.Value = "'  -0.66%  "
$c.Style = 'Normal'
$c = $ws.Range('D33')
$c.Value = "'0.7490"
$c.Style = 'Normal'
$c = $ws.Range('E33')
$c.Value = "'  -2.94%  "
$c.Style = 'Normal'
$c = $ws.Range('E34')
$c.Value = "'  -0.47%  "
$c.Style = 'Normal'
$c = $ws.Range('D35')
$c.Value = "'1.129"
$c.Style = 'Normal'
$c = $ws.Range('E35')
$c.Value = "'  -3.72%  "
$c.Style = 'Normal'
$c = $ws.Range('D36')
$c.Value = "'2.537"
$c.Style = 'Normal'
$c = $ws.Range('E36')
$c.Value = "'  -3.55%  "
$c.Style = 'Normal'
$c = $ws.Range('D37')
$c.Value = "'1.083"
$c.Style = 'Normal'
$c = $ws.Range('E37')
$c.Value = "'  +0.52%  "
$c.Style = 'Normal'
$c = $ws.Range('D38')
$c.Value = "'0.01936"
$c.Style = 'Normal'
$c = $ws.Range('E38')
$c.Value = "'  -1.17%  "
$c.Style = 'Normal'
$c = $ws.Range('D39')
$c.Value = "'2.908"
$c.Style = 'Normal'
$c = $ws.Range('E39')
$c.Value = "'  -1.43%  "
$c.Style = 'Normal'
$c = $ws.Range('D40')
$c.Value = "'0.05111"
$c.Style = 'Normal'
$c = $ws.Range('E40')
$c.Value = "'  -2.31%  "
$c.Style = 'Normal'
$c = $ws.Range('D41')
$c.Value = "'6.927"
$c.Style = 'Normal'
$c = $ws.Range('E41')
$c.Value = "'  -1.72%  "
$c.Style = 'Normal'
$c = $ws.Range('D42')
$c.Value = "'0.4967"
$c.Style = 'Normal'
$c = $ws.Range('E42')
$c.Value = "'  -3.37%  "
$c.Style = 'Normal'
$c = $ws.Range('E43')
$c.Value = "'  -2.79%  "
$c.Style = 'Normal'
$c = $ws.Range('D44')
$c.Value = "'8.283"
$c.Style = 'Normal'
$c = $ws.Range('E44')
$c.Value = "'  -1.34%  "
$c.Style = 'Normal'
$c = $ws.Range('D45')
$c.Value = "'0.4680"
$c.Style = 'Normal'
$c = $ws.Range('E45')
$c.Value = "'  -3.12%  "
$c.Style = 'Normal'
$c = $ws.Range('D46')
$c.Value = "'1.006"
$c.Style = 'Normal'
$c = $ws.Range('E46')
$c.Value = "'  +0.56%  "
$c.Style = 'Normal'
$c = $ws.Range('D47')
$c.Value = "'10.07"
$c.Style = 'Normal'
$c = $ws.Range('E47')
$c.Value = "'  -2.69%  "
$c.Style = 'Normal'
$c = $ws.Range('D48')
$c.Value = "'101.28"
$c.Style = 'Normal'
$c = $ws.Range('E48')
$c.Value = "'  -1.95%  "
$c.Style = 'Normal'
$c = $ws.Range('D49')
$c.Value = "'1.609"
$c.Style = 'Normal'
$c = $ws.Range('E49')
$c.Value = "'  -2.58%  "
$c.Style = 'Normal'
$c = $ws.Range('D50')
$c.Value = "'0.06089"
$c.Style = 'Normal'
$c = $ws.Range('E50')
$c.Value = "'  -2.09%  "
$c.Style = 'Normal'
$c = $ws.Range('D51')
$c.Value = "'64.38"
$c.Style = 'Normal'
$c = $ws.Range('E51')
$c.Value = "'  -2.41%  "
$c.Style = 'Normal'
